$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.267632722854614
$ws.Range("B1").Value = 2.922163486480713
$ws.Range("C1").Value = 5.519822597503662
$ws.Range("D1").Value = 1.868225336074829
$ws.Range("E1").Value = 1.030337691307068
